# Generate Report for Handback
# The "6db3fb88-e941-4d39-93c1-61e314267448" localization item has been
# handed back (in sync with en-US) for both zh-cn and de-de locales.
# Update the Overview sheet and each locale sheet accordingly.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $status
$ovw.Range("C2").Value = $status

# --- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $status
$zh.Range("F2").Value = "6db3fb88-e941-4d39-93c1-61e314267448.md"
$zh.Range("G2").Value = "6db3fb88-e941-4d39-93c1-61e314267448.c70940a3c4a7b3239228265732dbb2e629cb5224.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-19 04:31:56"

# --- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $status
$de.Range("F2").Value = "6db3fb88-e941-4d39-93c1-61e314267448.md"
$de.Range("G2").Value = "6db3fb88-e941-4d39-93c1-61e314267448.c70940a3c4a7b3239228265732dbb2e629cb5224.de-de.xlf"
$de.Range("H2").Value = "2016-03-19 04:32:00"

Write-Output "Handback report generated"
